# Apply updated cryptocurrency data (prices and 1h volume change %) to sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.711.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +3.00%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.447.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.87%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.09%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'578.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +2.75%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'145.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +2.85%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.08%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.55%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.446.40"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.49%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +2.50%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +1.40%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +1.29%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +3.04%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'28.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +8.56%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  +5.89%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'2.890.89"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'62.543.98"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +3.52%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.441.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.52%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'7.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -4.27%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'10.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +2.69%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'328.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.37%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +1.13%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +10.43%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +0.01%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'65.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.08%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'644.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +12.97%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'1.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +16.48%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("B28").Value = "'Aptos"
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Value = "'8.47"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +4.89%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("B29").Value = "'BabyDogeCoin"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = "'0.0₆0537"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +92.87%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.0₃0989"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +4.96%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +2.23%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +1.41%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +7.77%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.88"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +3.65%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +5.38%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +2.14%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +0.22%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +3.44%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'5.49"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +6.13%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'153.37"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.03%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +0.86%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'18.65"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.75%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'2.73"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +6.49%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +5.35%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D46").Value = "'42.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.64%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  +27.86%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'145.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +2.15%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +2.26%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +6.95%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.606"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +2.93%  "
$ws.Range("E51").Style = "Normal"
